# Ccl11-Ackr4.xlsx: refresh NATMI ligand-receptor metrics (columns E:T, data rows 2:19)
# with newly recomputed TPM-based values ("update scripts wuth new tpm").
# Columns A:D (Sending cluster / Ligand symbol / Receptor symbol / Target cluster) are
# unchanged identifiers, so only E:T are rewritten, row by row, via a single Range write.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 18,16
# row 2: ECs -> ECs
$data[0,0] = 3
$data[0,1] = 1
$data[0,2] = 2.940931333333333
$data[0,3] = 8.822794
$data[0,4] = 0.03454096854573427
$data[0,5] = 0.03454096854573427
$data[0,6] = 2
$data[0,7] = 0.6666666666666666
$data[0,8] = 0.180428
$data[0,9] = 0.541284
$data[0,10] = 0.6724860231084607
$data[0,11] = 0.6724860231084607
$data[0,12] = 0.5306263586106666
$data[0,13] = 4.775637227496
$data[0,14] = 0.02322831857163527
$data[0,15] = 0.02322831857163527
# row 3: ECs -> FAPs
$data[1,0] = 3
$data[1,1] = 1
$data[1,2] = 2.940931333333333
$data[1,3] = 8.822794
$data[1,4] = 0.03454096854573427
$data[1,5] = 0.03454096854573427
$data[1,6] = 1
$data[1,7] = 0.3333333333333333
$data[1,8] = 0.01727566666666666
$data[1,9] = 0.051827
$data[1,10] = 0.06438936513852653
$data[1,11] = 0.06438936513852653
$data[1,12] = 0.05080654940422222
$data[1,13] = 0.457258944638
$data[1,14] = 0.002224071035929644
$data[1,15] = 0.002224071035929644
# row 4: ECs -> MuSCs
$data[2,0] = 3
$data[2,1] = 1
$data[2,2] = 2.940931333333333
$data[2,3] = 8.822794
$data[2,4] = 0.03454096854573427
$data[2,5] = 0.03454096854573427
$data[2,6] = 2
$data[2,7] = 0.6666666666666666
$data[2,8] = 0.07059633333333333
$data[2,9] = 0.211789
$data[2,10] = 0.2631246117530128
$data[2,11] = 0.2631246117530128
$data[2,12] = 0.2076189687184444
$data[2,13] = 1.868570718466
$data[2,14] = 0.009088578938169358
$data[2,15] = 0.009088578938169358
# row 5: FAPs -> ECs
$data[3,0] = 3
$data[3,1] = 1
$data[3,2] = 76.49331166666667
$data[3,3] = 229.479935
$data[3,4] = 0.8984069237831173
$data[3,5] = 0.8984069237831174
$data[3,6] = 2
$data[3,7] = 0.6666666666666666
$data[3,8] = 0.180428
$data[3,9] = 0.541284
$data[3,10] = 0.6724860231084607
$data[3,11] = 0.6724860231084607
$data[3,12] = 13.80153523739333
$data[3,13] = 124.21381713654
$data[3,14] = 0.6041660993080145
$data[3,15] = 0.6041660993080146
# row 6: FAPs -> FAPs
$data[4,0] = 3
$data[4,1] = 1
$data[4,2] = 76.49331166666667
$data[4,3] = 229.479935
$data[4,4] = 0.8984069237831173
$data[4,5] = 0.8984069237831174
$data[4,6] = 1
$data[4,7] = 0.3333333333333333
$data[4,8] = 0.01727566666666666
$data[4,9] = 0.051827
$data[4,10] = 0.06438936513852653
$data[4,11] = 0.06438936513852653
$data[4,12] = 1.321472954582778
$data[4,13] = 11.893256591245
$data[4,14] = 0.05784785145845151
$data[4,15] = 0.05784785145845152
# row 7: FAPs -> MuSCs
$data[5,0] = 3
$data[5,1] = 1
$data[5,2] = 76.49331166666667
$data[5,3] = 229.479935
$data[5,4] = 0.8984069237831173
$data[5,5] = 0.8984069237831174
$data[5,6] = 2
$data[5,7] = 0.6666666666666666
$data[5,8] = 0.07059633333333333
$data[5,9] = 0.211789
$data[5,10] = 0.2631246117530128
$data[5,11] = 0.2631246117530128
$data[5,12] = 5.400147328190555
$data[5,13] = 48.601325953715
$data[5,14] = 0.2363929730166513
$data[5,15] = 0.2363929730166514
# row 8: Inflammatory-Mac -> ECs
$data[6,0] = 3
$data[6,1] = 1
$data[6,2] = 1.002166333333333
$data[6,3] = 3.006499
$data[6,4] = 0.01177035159063915
$data[6,5] = 0.01177035159063915
$data[6,6] = 2
$data[6,7] = 0.6666666666666666
$data[6,8] = 0.180428
$data[6,9] = 0.541284
$data[6,10] = 0.6724860231084607
$data[6,11] = 0.6724860231084607
$data[6,12] = 0.1808188671906667
$data[6,13] = 1.627369804716
$data[6,14] = 0.007915396931777265
$data[6,15] = 0.007915396931777266
# row 9: Inflammatory-Mac -> FAPs
$data[7,0] = 3
$data[7,1] = 1
$data[7,2] = 1.002166333333333
$data[7,3] = 3.006499
$data[7,4] = 0.01177035159063915
$data[7,5] = 0.01177035159063915
$data[7,6] = 1
$data[7,7] = 0.3333333333333333
$data[7,8] = 0.01727566666666666
$data[7,9] = 0.051827
$data[7,10] = 0.06438936513852653
$data[7,11] = 0.06438936513852653
$data[7,12] = 0.01731309151922222
$data[7,13] = 0.155817823673
$data[7,14] = 0.0007578854663785005
$data[7,15] = 0.0007578854663785006
# row 10: Inflammatory-Mac -> MuSCs
$data[8,0] = 3
$data[8,1] = 1
$data[8,2] = 1.002166333333333
$data[8,3] = 3.006499
$data[8,4] = 0.01177035159063915
$data[8,5] = 0.01177035159063915
$data[8,6] = 2
$data[8,7] = 0.6666666666666666
$data[8,8] = 0.07059633333333333
$data[8,9] = 0.211789
$data[8,10] = 0.2631246117530128
$data[8,11] = 0.2631246117530128
$data[8,12] = 0.07074926852344443
$data[8,13] = 0.636743416711
$data[8,14] = 0.003097069192483383
$data[8,15] = 0.003097069192483383
# row 11: MuSCs -> ECs
$data[9,0] = 3
$data[9,1] = 1
$data[9,2] = 2.356521666666667
$data[9,3] = 7.069565
$data[9,4] = 0.02767713065691252
$data[9,5] = 0.02767713065691253
$data[9,6] = 2
$data[9,7] = 0.6666666666666666
$data[9,8] = 0.180428
$data[9,9] = 0.541284
$data[9,10] = 0.6724860231084607
$data[9,11] = 0.6724860231084607
$data[9,12] = 0.4251824912733334
$data[9,13] = 3.82664242146
$data[9,14] = 0.01861248352652036
$data[9,15] = 0.01861248352652036
# row 12: MuSCs -> FAPs
$data[10,0] = 3
$data[10,1] = 1
$data[10,2] = 2.356521666666667
$data[10,3] = 7.069565
$data[10,4] = 0.02767713065691252
$data[10,5] = 0.02767713065691253
$data[10,6] = 1
$data[10,7] = 0.3333333333333333
$data[10,8] = 0.01727566666666666
$data[10,9] = 0.051827
$data[10,10] = 0.06438936513852653
$data[10,11] = 0.06438936513852653
$data[10,12] = 0.04071048280611111
$data[10,13] = 0.366394345255
$data[10,14] = 0.001782112871854647
$data[10,15] = 0.001782112871854647
# row 13: MuSCs -> MuSCs
$data[11,0] = 3
$data[11,1] = 1
$data[11,2] = 2.356521666666667
$data[11,3] = 7.069565
$data[11,4] = 0.02767713065691252
$data[11,5] = 0.02767713065691253
$data[11,6] = 2
$data[11,7] = 0.6666666666666666
$data[11,8] = 0.07059633333333333
$data[11,9] = 0.211789
$data[11,10] = 0.2631246117530128
$data[11,11] = 0.2631246117530128
$data[11,12] = 0.1663617890872222
$data[11,13] = 1.497256101785
$data[11,14] = 0.007282534258537516
$data[11,15] = 0.007282534258537517
# row 14: Neutrophils -> ECs
$data[12,0] = 2
$data[12,1] = 0.6666666666666666
$data[12,2] = 0.5156633333333334
$data[12,3] = 1.54699
$data[12,4] = 0.006056418514425867
$data[12,5] = 0.006056418514425868
$data[12,6] = 2
$data[12,7] = 0.6666666666666666
$data[12,8] = 0.180428
$data[12,9] = 0.541284
$data[12,10] = 0.6724860231084607
$data[12,11] = 0.6724860231084607
$data[12,12] = 0.09304010390666667
$data[12,13] = 0.83736093516
$data[12,14] = 0.004072856801046703
$data[12,15] = 0.004072856801046703
# row 15: Neutrophils -> FAPs
$data[13,0] = 2
$data[13,1] = 0.6666666666666666
$data[13,2] = 0.5156633333333334
$data[13,3] = 1.54699
$data[13,4] = 0.006056418514425867
$data[13,5] = 0.006056418514425868
$data[13,6] = 1
$data[13,7] = 0.3333333333333333
$data[13,8] = 0.01727566666666666
$data[13,9] = 0.051827
$data[13,10] = 0.06438936513852653
$data[13,11] = 0.06438936513852653
$data[13,12] = 0.008908427858888889
$data[13,13] = 0.08017585073
$data[13,14] = 0.0003899689431570995
$data[13,15] = 0.0003899689431570996
# row 16: Neutrophils -> MuSCs
$data[14,0] = 2
$data[14,1] = 0.6666666666666666
$data[14,2] = 0.5156633333333334
$data[14,3] = 1.54699
$data[14,4] = 0.006056418514425867
$data[14,5] = 0.006056418514425868
$data[14,6] = 2
$data[14,7] = 0.6666666666666666
$data[14,8] = 0.07059633333333333
$data[14,9] = 0.211789
$data[14,10] = 0.2631246117530128
$data[14,11] = 0.2631246117530128
$data[14,12] = 0.03640394056777778
$data[14,13] = 0.32763546511
$data[14,14] = 0.001593592770222065
$data[14,15] = 0.001593592770222065
# row 17: Resolving-Mac -> ECs
$data[15,0] = 3
$data[15,1] = 1
$data[15,2] = 1.834685
$data[15,3] = 5.504055
$data[15,4] = 0.02154820690917088
$data[15,5] = 0.02154820690917089
$data[15,6] = 2
$data[15,7] = 0.6666666666666666
$data[15,8] = 0.180428
$data[15,9] = 0.541284
$data[15,10] = 0.6724860231084607
$data[15,11] = 0.6724860231084607
$data[15,12] = 0.33102854518
$data[15,13] = 2.97925690662
$data[15,14] = 0.01449086796946658
$data[15,15] = 0.01449086796946659
# row 18: Resolving-Mac -> FAPs
$data[16,0] = 3
$data[16,1] = 1
$data[16,2] = 1.834685
$data[16,3] = 5.504055
$data[16,4] = 0.02154820690917088
$data[16,5] = 0.02154820690917089
$data[16,6] = 1
$data[16,7] = 0.3333333333333333
$data[16,8] = 0.01727566666666666
$data[16,9] = 0.051827
$data[16,10] = 0.06438936513852653
$data[16,11] = 0.06438936513852653
$data[16,12] = 0.03169540649833334
$data[16,13] = 0.285258658485
$data[16,14] = 0.001387475362755124
$data[16,15] = 0.001387475362755125
# row 19: Resolving-Mac -> MuSCs
$data[17,0] = 3
$data[17,1] = 1
$data[17,2] = 1.834685
$data[17,3] = 5.504055
$data[17,4] = 0.02154820690917088
$data[17,5] = 0.02154820690917089
$data[17,6] = 2
$data[17,7] = 0.6666666666666666
$data[17,8] = 0.07059633333333333
$data[17,9] = 0.211789
$data[17,10] = 0.2631246117530128
$data[17,11] = 0.2631246117530128
$data[17,12] = 0.1295220338216667
$data[17,13] = 1.165698304395
$data[17,14] = 0.005669863576949177
$data[17,15] = 0.005669863576949178

$ws.Range("E2:T19").Value = $data
